$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 23 - this shifts the existing rows 23:33
# down to 24:34, matching the diff (dimension grows from A1:R33 to A1:R34).
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new data record.
# (All other columns mirror the record that used to be in row 23 before
# the shift - only the date (D) and volume (J) differ for the new entry.)
$ws.Cells.Item(23, 1).Value = 11
$ws.Cells.Item(23, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(23, 3).Value = "Bíobío"
$ws.Cells.Item(23, 4).Value = 44755
$ws.Cells.Item(23, 5).Value = 8
$ws.Cells.Item(23, 6).Value = 100114007
$ws.Cells.Item(23, 7).Value = "Jengibre"
$ws.Cells.Item(23, 8).Value = "Sin especificar"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 40
$ws.Cells.Item(23, 11).Value = 14000
$ws.Cells.Item(23, 12).Value = 15000
$ws.Cells.Item(23, 13).Value = 14600
$ws.Cells.Item(23, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(23, 15).Value = "Perú"
$ws.Cells.Item(23, 16).Value = 1123
$ws.Cells.Item(23, 17).Value = 13
$ws.Cells.Item(23, 18).Value = "Hortaliza"
